$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -7.548999999999999
$ws.Range("C3").Value = -12.31519999999999
$ws.Range("E3").Value = 15.57460000000001
$ws.Range("E12").Value = 17.20000000000002
$ws.Range("C14").Value = -12.7396
$ws.Range("C21").Value = -12.53120000000001
$ws.Range("C23").Value = -12.5253
$ws.Range("E24").Value = 16.7894
$ws.Range("C25").Value = -12.5858
$ws.Range("D25").Value = -7.902400000000004
$ws.Range("E25").Value = 16.85910000000001
$ws.Range("C26").Value = -12.52390000000001
$ws.Range("D27").Value = -8.787300000000009
$ws.Range("C29").Value = -10.83050000000001
$ws.Range("D31").Value = -8.68260000000001
$ws.Range("D39").Value = -8.113199999999997
$ws.Range("D48").Value = -7.224599999999997
$ws.Range("E50").Value = 16.3473
$ws.Range("D51").Value = -7.574699999999998
$ws.Range("D52").Value = -7.7537
$ws.Range("C53").Value = -11.0085
$ws.Range("E53").Value = 16.87920000000002
$ws.Range("D55").Value = -8.387999999999998
$ws.Range("D56").Value = -7.865400000000002
$ws.Range("C57").Value = -14.07139999999999
$ws.Range("D57").Value = -7.7451
$ws.Range("E57").Value = 16.8723
$ws.Range("C59").Value = -12.7868
$ws.Range("E61").Value = 16.53340000000001
$ws.Range("E63").Value = 17.37390000000002
$ws.Range("C69").Value = -10.8483
$ws.Range("E70").Value = 17.30340000000001
$ws.Range("D73").Value = -7.812300000000001
$ws.Range("C79").Value = -10.96610000000001
$ws.Range("C83").Value = -13.91499999999999
$ws.Range("E86").Value = 16.70229999999999
$ws.Range("D89").Value = -5.679400000000004
$ws.Range("D90").Value = -8.308600000000002
$ws.Range("C91").Value = -10.2953
$ws.Range("D92").Value = -5.878700000000002
$ws.Range("C93").Value = -11.40020000000001
$ws.Range("E98").Value = 15.2891
$ws.Range("E100").Value = 16.96080000000001
$ws.Range("E102").Value = 16.64629999999998